# Generate Report for Handoff
#
# The localization-status report records, per source file, the timestamp of
# the most recent handoff. Re-running the handoff-status generator for
# "c86cf631-726d-415a-9e5e-3deb6a4488fa" (row 6 on every sheet) refreshed its
# "Latest Handoff Date" / "Latest Handoff Datetime" to a newer timestamp on
# the Overview sheet and on each per-locale sheet (zh-cn, de-de). No other
# cell's displayed content changes.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" column (D) for the c86cf631... row.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-40-13 14:40:48"

# zh-cn sheet: "Latest Handoff Datetime" column (E) for the c86cf631... row.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-13 14:40:45"

# de-de sheet: "Latest Handoff Datetime" column (E) for the c86cf631... row.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-13 14:40:48"
